$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = "Gestão de Fornecedores"
$ws.Range("B25").Value = "Gestão de Aditivos de Contratos"
$ws.Range("B32").Value = "Gestão de Itens de Auditoria"
$ws.Range("B37").Value = "Gestão de Bilhetes Telefônicos"
$ws.Range("B39").Value = "Gestão de Estoque de Aparelhos"
$ws.Range("B40").Value = "Gestão de Notas Fiscais Estoque"
$ws.Range("B41").Value = "Gestão de Endereços de Entrega"
$ws.Range("B45").Value = "Gestão de Tipos de Consumidores"
$ws.Range("B48").Value = "Gestão de Linhas Móveis e Chips SIM"
$ws.Range("B52").Value = "Gestão de Lotes de Auditoria"
$ws.Range("B57").Value = "Gestão de Status e Tipos Genéricos"
$ws.Range("B58").Value = "Gestão de Tipos de Chamado"
$ws.Range("B64").Value = "Notificações e Alertas"
$ws.Range("B66").Value = "Inventário Cíclico e Auditoria de Estoque"
$ws.Range("B69").Value = "Pesquisa de Satisfação"
$ws.Range("B72").Value = "Gestão de Chamados - Portal Self-Service"
$ws.Range("B81").Value = "Importação de Dados"
$ws.Range("B84").Value = "Aprovações e Workflows"
$ws.Range("B90").Value = "Auditoria de Custos e Compliance"
$ws.Range("B96").Value = "Dashboards e KPIs com Análise Preditiva"
$ws.Range("B99").Value = "Relatórios e Volumetria"
$ws.Range("B102").Value = "Marcadores Localização QRCode"
$ws.Range("B103").Value = "CAPTCHA, MFA, Contestação e Segurança Avançada"
$ws.Range("B104").Value = "Gestão de Documentos Originais e Digitalização"
$ws.Range("B106").Value = "Backup, Recuperação e Disaster Recovery"
